$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 125
$ws.Cells.Item(125, 1).Value = "North Carolina"
$ws.Cells.Item(125, 2).Value = 122.472408880856
$ws.Cells.Item(125, 3).Value = 94.6831103934111
$ws.Cells.Item(125, 4).Value = 0.731504316498081
$ws.Cells.Item(125, 5).Value = 0.603121914343077
$ws.Cells.Item(125, 6).Value = 0.790876256090942
$ws.Cells.Item(125, 7).Value = 110.797185112423
$ws.Cells.Item(125, 8).Value = 99.2863885840249
$ws.Cells.Item(125, 9).Value = 71.7572758629269
$ws.Cells.Item(125, 10).Value = 13.1
$ws.Cells.Item(125, 11).Value = 55.8
$ws.Cells.Item(125, 12).Value = 61.1
$ws.Cells.Item(125, 13).Value = 19.9
$ws.Cells.Item(125, 14).Value = 19
$ws.Cells.Item(125, 15).Value = 52.5
$ws.Cells.Item(125, 16).Value = 32.1
$ws.Cells.Item(125, 17).Value = 54.1
$ws.Cells.Item(125, 18).Value = 74.7
$ws.Cells.Item(125, 19).Value = 112.1
$ws.Cells.Item(125, 20).Value = 64
$ws.Cells.Item(125, 21).Value = 21.1
$ws.Cells.Item(125, 22).Value = 39.7
$ws.Cells.Item(125, 23).Value = 73.4
$ws.Cells.Item(125, 24).Value = 55.7
$ws.Cells.Item(125, 25).Value = 7.4
$ws.Cells.Item(125, 26).Value = 8.4
$ws.Cells.Item(125, 27).Value = 1.656
$ws.Cells.Item(125, 28).Value = 0.575
$ws.Cells.Item(125, 29).Value = 12.6
$ws.Cells.Item(125, 30).Value = 22.2
$ws.Cells.Item(125, 31).Value = -13.1
$ws.Cells.Item(125, 32).Value = 45.2
$ws.Cells.Item(125, 33).Value = 48.6
$ws.Cells.Item(125, 34).Value = 33
$ws.Cells.Item(125, 35).Value = 18.4
$ws.Cells.Item(125, 36).Value = 47.8
$ws.Cells.Item(125, 37).Value = 36
$ws.Cells.Item(125, 38).Value = 44.3
$ws.Cells.Item(125, 39).Value = 70.2
$ws.Cells.Item(125, 40).Value = 102.2
$ws.Cells.Item(125, 41).Value = 59.7
$ws.Cells.Item(125, 42).Value = 18.3
$ws.Cells.Item(125, 43).Value = 26.6
$ws.Cells.Item(125, 44).Value = 60.3
$ws.Cells.Item(125, 45).Value = 4.6
$ws.Cells.Item(125, 46).Value = 6.6
$ws.Cells.Item(125, 47).Value = 0.863
$ws.Cells.Item(125, 48).Value = 0.451
$ws.Cells.Item(125, 49).Value = 15.9
$ws.Cells.Item(125, 50).Value = 21.8
$ws.Cells.Item(125, 51).Value = 5.6
$ws.Cells.Item(125, 52).Value = 1.027
$ws.Cells.Item(125, 53).Value = 0.949
$ws.Cells.Item(125, 54).Value = 0.5
$ws.Cells.Item(125, 55).Value = 0.5
$ws.Cells.Item(125, 56).Value = "Syracuse"
$ws.Cells.Item(125, 57).Value = 111.887498470978
$ws.Cells.Item(125, 58).Value = 93.2967229009327
$ws.Cells.Item(125, 59).Value = 0.700702559833434
$ws.Cells.Item(125, 60).Value = 0.517936613792802
$ws.Cells.Item(125, 61).Value = 0.794525575768337
$ws.Cells.Item(125, 62).Value = 110.052235729994
$ws.Cells.Item(125, 63).Value = 99.452497067608
$ws.Cells.Item(125, 64).Value = 65.4975872773356
$ws.Cells.Item(125, 65).Value = 5.5
$ws.Cells.Item(125, 66).Value = 48.3
$ws.Cells.Item(125, 67).Value = 43.6
$ws.Cells.Item(125, 68).Value = 36.5
$ws.Cells.Item(125, 69).Value = 19.9
$ws.Cells.Item(125, 70).Value = 50.2
$ws.Cells.Item(125, 71).Value = 36.1
$ws.Cells.Item(125, 72).Value = 47.3
$ws.Cells.Item(125, 73).Value = 69.4
$ws.Cells.Item(125, 74).Value = 107
$ws.Cells.Item(125, 75).Value = 55.9
$ws.Cells.Item(125, 76).Value = 20.1
$ws.Cells.Item(125, 77).Value = 31.7
$ws.Cells.Item(125, 78).Value = 67
$ws.Cells.Item(125, 79).Value = 49.5
$ws.Cells.Item(125, 80).Value = 7.7
$ws.Cells.Item(125, 81).Value = 10.4
$ws.Cells.Item(125, 82).Value = 1.143
$ws.Cells.Item(125, 83).Value = 0.569
$ws.Cells.Item(125, 84).Value = 15.3
$ws.Cells.Item(125, 85).Value = 20.7
$ws.Cells.Item(125, 86).Value = -5.5
$ws.Cells.Item(125, 87).Value = 44.9
$ws.Cells.Item(125, 88).Value = 50.5
$ws.Cells.Item(125, 89).Value = 32.9
$ws.Cells.Item(125, 90).Value = 16.6
$ws.Cells.Item(125, 91).Value = 47
$ws.Cells.Item(125, 92).Value = 30.8
$ws.Cells.Item(125, 93).Value = 47.6
$ws.Cells.Item(125, 94).Value = 66.3
$ws.Cells.Item(125, 95).Value = 99.4
$ws.Cells.Item(125, 96).Value = 57.3
$ws.Cells.Item(125, 97).Value = 16.2
$ws.Cells.Item(125, 98).Value = 33
$ws.Cells.Item(125, 99).Value = 68.3
$ws.Cells.Item(125, 100).Value = 5.3
$ws.Cells.Item(125, 101).Value = 8.1
$ws.Cells.Item(125, 102).Value = 1.152
$ws.Cells.Item(125, 103).Value = 0.656
$ws.Cells.Item(125, 104).Value = 17
$ws.Cells.Item(125, 105).Value = 22.8
$ws.Cells.Item(125, 106).Value = 0.8
$ws.Cells.Item(125, 107).Value = 0.976
$ws.Cells.Item(125, 108).Value = 0.964
$ws.Cells.Item(125, 109).Value = 0.5
$ws.Cells.Item(125, 110).Value = 0.5
$ws.Cells.Item(125, 111).Value = "'FALSE"

# Row 126
$ws.Cells.Item(126, 1).Value = "Oklahoma"
$ws.Cells.Item(126, 2).Value = 118.70227059828
$ws.Cells.Item(126, 3).Value = 93.1606214962715
$ws.Cells.Item(126, 4).Value = 0.75453708018108
$ws.Cells.Item(126, 5).Value = 0.563394891283762
$ws.Cells.Item(126, 6).Value = 0.839881470981451
$ws.Cells.Item(126, 7).Value = 110.708500425401
$ws.Cells.Item(126, 8).Value = 97.7125081879537
$ws.Cells.Item(126, 9).Value = 70.901528263298
$ws.Cells.Item(126, 10).Value = 10.1
$ws.Cells.Item(126, 11).Value = 49.8
$ws.Cells.Item(126, 12).Value = 43
$ws.Cells.Item(126, 13).Value = 38.9
$ws.Cells.Item(126, 14).Value = 18.1
$ws.Cells.Item(126, 15).Value = 54.8
$ws.Cells.Item(126, 16).Value = 42.8
$ws.Cells.Item(126, 17).Value = 48.4
$ws.Cells.Item(126, 18).Value = 72.3
$ws.Cells.Item(126, 19).Value = 115.5
$ws.Cells.Item(126, 20).Value = 60.1
$ws.Cells.Item(126, 21).Value = 20.1
$ws.Cells.Item(126, 22).Value = 29.6
$ws.Cells.Item(126, 23).Value = 72
$ws.Cells.Item(126, 24).Value = 51.5
$ws.Cells.Item(126, 25).Value = 8.2
$ws.Cells.Item(126, 26).Value = 8
$ws.Cells.Item(126, 27).Value = 1.138
$ws.Cells.Item(126, 28).Value = 0.53
$ws.Cells.Item(126, 29).Value = 15.6
$ws.Cells.Item(126, 30).Value = 20.3
$ws.Cells.Item(126, 31).Value = -10.1
$ws.Cells.Item(126, 32).Value = 44.9
$ws.Cells.Item(126, 33).Value = 51.4
$ws.Cells.Item(126, 34).Value = 31.6
$ws.Cells.Item(126, 35).Value = 17
$ws.Cells.Item(126, 36).Value = 46.4
$ws.Cells.Item(126, 37).Value = 33.1
$ws.Cells.Item(126, 38).Value = 44.7
$ws.Cells.Item(126, 39).Value = 67.7
$ws.Cells.Item(126, 40).Value = 98.7
$ws.Cells.Item(126, 41).Value = 62.9
$ws.Cells.Item(126, 42).Value = 17.7
$ws.Cells.Item(126, 43).Value = 28
$ws.Cells.Item(126, 44).Value = 70.4
$ws.Cells.Item(126, 45).Value = 5.3
$ws.Cells.Item(126, 46).Value = 8.6
$ws.Cells.Item(126, 47).Value = 0.869
$ws.Cells.Item(126, 48).Value = 0.432
$ws.Cells.Item(126, 49).Value = 15.2
$ws.Cells.Item(126, 50).Value = 22.9
$ws.Cells.Item(126, 51).Value = -0.9
$ws.Cells.Item(126, 52).Value = 0.953
$ws.Cells.Item(126, 53).Value = 0.966
$ws.Cells.Item(126, 54).Value = 0.643
$ws.Cells.Item(126, 55).Value = 0.357
$ws.Cells.Item(126, 56).Value = "Villanova"
$ws.Cells.Item(126, 57).Value = 121.607651393701
$ws.Cells.Item(126, 58).Value = 90.9489303078528
$ws.Cells.Item(126, 59).Value = 0.711354159630472
$ws.Cells.Item(126, 60).Value = 0.588242734646016
$ws.Cells.Item(126, 61).Value = 0.752598454548143
$ws.Cells.Item(126, 62).Value = 109.45666079805
$ws.Cells.Item(126, 63).Value = 98.7883640982027
$ws.Cells.Item(126, 64).Value = 66.8624144413201
$ws.Cells.Item(126, 65).Value = 13.9
$ws.Cells.Item(126, 66).Value = 51.6
$ws.Cells.Item(126, 67).Value = 47.2
$ws.Cells.Item(126, 68).Value = 33.4
$ws.Cells.Item(126, 69).Value = 19.4
$ws.Cells.Item(126, 70).Value = 55.2
$ws.Cells.Item(126, 71).Value = 35.4
$ws.Cells.Item(126, 72).Value = 56.8
$ws.Cells.Item(126, 73).Value = 78.4
$ws.Cells.Item(126, 74).Value = 118
$ws.Cells.Item(126, 75).Value = 56.6
$ws.Cells.Item(126, 76).Value = 19.2
$ws.Cells.Item(126, 77).Value = 28.1
$ws.Cells.Item(126, 78).Value = 74.5
$ws.Cells.Item(126, 79).Value = 51.4
$ws.Cells.Item(126, 80).Value = 6
$ws.Cells.Item(126, 81).Value = 9
$ws.Cells.Item(126, 82).Value = 1.476
$ws.Cells.Item(126, 83).Value = 0.603
$ws.Cells.Item(126, 84).Value = 14.3
$ws.Cells.Item(126, 85).Value = 22.2
$ws.Cells.Item(126, 86).Value = -13.9
$ws.Cells.Item(126, 87).Value = 43.2
$ws.Cells.Item(126, 88).Value = 49
$ws.Cells.Item(126, 89).Value = 33.3
$ws.Cells.Item(126, 90).Value = 17.8
$ws.Cells.Item(126, 91).Value = 46.8
$ws.Cells.Item(126, 92).Value = 33.7
$ws.Cells.Item(126, 93).Value = 44.6
$ws.Cells.Item(126, 94).Value = 66.2
$ws.Cells.Item(126, 95).Value = 99.4
$ws.Cells.Item(126, 96).Value = 55.9
$ws.Cells.Item(126, 97).Value = 17.1
$ws.Cells.Item(126, 98).Value = 25.5
$ws.Cells.Item(126, 99).Value = 71.9
$ws.Cells.Item(126, 100).Value = 3.8
$ws.Cells.Item(126, 101).Value = 7.5
$ws.Cells.Item(126, 102).Value = 0.927
$ws.Cells.Item(126, 103).Value = 0.565
$ws.Cells.Item(126, 104).Value = 17.8
$ws.Cells.Item(126, 105).Value = 24.5
$ws.Cells.Item(126, 106).Value = 2.4
$ws.Cells.Item(126, 107).Value = 0.958
$ws.Cells.Item(126, 108).Value = 0.923
$ws.Cells.Item(126, 109).Value = 0.8
$ws.Cells.Item(126, 110).Value = 0.2
$ws.Cells.Item(126, 111).Value = "'FALSE"

# Row 127
$ws.Cells.Item(127, 1).Value = "North Carolina"
$ws.Cells.Item(127, 2).Value = 122.95885823211
$ws.Cells.Item(127, 3).Value = 94.6707539780454
$ws.Cells.Item(127, 4).Value = 0.735613786620976
$ws.Cells.Item(127, 5).Value = 0.603010895211314
$ws.Cells.Item(127, 6).Value = 0.790845910588888
$ws.Cells.Item(127, 7).Value = 110.827153441749
$ws.Cells.Item(127, 8).Value = 99.1406160590647
$ws.Cells.Item(127, 9).Value = 71.8637683642112
$ws.Cells.Item(127, 10).Value = 13.2
$ws.Cells.Item(127, 11).Value = 55.9
$ws.Cells.Item(127, 12).Value = 61.5
$ws.Cells.Item(127, 13).Value = 19.7
$ws.Cells.Item(127, 14).Value = 18.8
$ws.Cells.Item(127, 15).Value = 52.6
$ws.Cells.Item(127, 16).Value = 31.9
$ws.Cells.Item(127, 17).Value = 54.4
$ws.Cells.Item(127, 18).Value = 74.8
$ws.Cells.Item(127, 19).Value = 112.3
$ws.Cells.Item(127, 20).Value = 64
$ws.Cells.Item(127, 21).Value = 20.8
$ws.Cells.Item(127, 22).Value = 40
$ws.Cells.Item(127, 23).Value = 73.1
$ws.Cells.Item(127, 24).Value = 55.8
$ws.Cells.Item(127, 25).Value = 7.4
$ws.Cells.Item(127, 26).Value = 8.3
$ws.Cells.Item(127, 27).Value = 1.652
$ws.Cells.Item(127, 28).Value = 0.573
$ws.Cells.Item(127, 29).Value = 12.7
$ws.Cells.Item(127, 30).Value = 22
$ws.Cells.Item(127, 31).Value = -13.2
$ws.Cells.Item(127, 32).Value = 45.3
$ws.Cells.Item(127, 33).Value = 48.8
$ws.Cells.Item(127, 34).Value = 33
$ws.Cells.Item(127, 35).Value = 18.1
$ws.Cells.Item(127, 36).Value = 47.7
$ws.Cells.Item(127, 37).Value = 35.9
$ws.Cells.Item(127, 38).Value = 44.3
$ws.Cells.Item(127, 39).Value = 69.5
$ws.Cells.Item(127, 40).Value = 101.9
$ws.Cells.Item(127, 41).Value = 59.9
$ws.Cells.Item(127, 42).Value = 18.2
$ws.Cells.Item(127, 43).Value = 26.9
$ws.Cells.Item(127, 44).Value = 60
$ws.Cells.Item(127, 45).Value = 4.7
$ws.Cells.Item(127, 46).Value = 6.7
$ws.Cells.Item(127, 47).Value = 0.856
$ws.Cells.Item(127, 48).Value = 0.444
$ws.Cells.Item(127, 49).Value = 15.8
$ws.Cells.Item(127, 50).Value = 21.7
$ws.Cells.Item(127, 51).Value = 5.4
$ws.Cells.Item(127, 52).Value = 1.028
$ws.Cells.Item(127, 53).Value = 0.952
$ws.Cells.Item(127, 54).Value = 0.5
$ws.Cells.Item(127, 55).Value = 0.5
$ws.Cells.Item(127, 56).Value = "Villanova"
$ws.Cells.Item(127, 57).Value = 122.783517287024
$ws.Cells.Item(127, 58).Value = 90.6258236924989
$ws.Cells.Item(127, 59).Value = 0.7193970466841
$ws.Cells.Item(127, 60).Value = 0.588201540444193
$ws.Cells.Item(127, 61).Value = 0.756867952390836
$ws.Cells.Item(127, 62).Value = 109.739975952838
$ws.Cells.Item(127, 63).Value = 98.6081868217284
$ws.Cells.Item(127, 64).Value = 66.8585044477694
$ws.Cells.Item(127, 65).Value = 14.7
$ws.Cells.Item(127, 66).Value = 52
$ws.Cells.Item(127, 67).Value = 47.3
$ws.Cells.Item(127, 68).Value = 33.4
$ws.Cells.Item(127, 69).Value = 19.3
$ws.Cells.Item(127, 70).Value = 55.8
$ws.Cells.Item(127, 71).Value = 35.9
$ws.Cells.Item(127, 72).Value = 57.3
$ws.Cells.Item(127, 73).Value = 78.2
$ws.Cells.Item(127, 74).Value = 119
$ws.Cells.Item(127, 75).Value = 56.4
$ws.Cells.Item(127, 76).Value = 19.2
$ws.Cells.Item(127, 77).Value = 28.3
$ws.Cells.Item(127, 78).Value = 74.1
$ws.Cells.Item(127, 79).Value = 51.4
$ws.Cells.Item(127, 80).Value = 6
$ws.Cells.Item(127, 81).Value = 9.2
$ws.Cells.Item(127, 82).Value = 1.471
$ws.Cells.Item(127, 83).Value = 0.596
$ws.Cells.Item(127, 84).Value = 14.3
$ws.Cells.Item(127, 85).Value = 22
$ws.Cells.Item(127, 86).Value = -14.7
$ws.Cells.Item(127, 87).Value = 43.1
$ws.Cells.Item(127, 88).Value = 49
$ws.Cells.Item(127, 89).Value = 33.3
$ws.Cells.Item(127, 90).Value = 17.7
$ws.Cells.Item(127, 91).Value = 46.5
$ws.Cells.Item(127, 92).Value = 33.3
$ws.Cells.Item(127, 93).Value = 44.5
$ws.Cells.Item(127, 94).Value = 66.1
$ws.Cells.Item(127, 95).Value = 98.8
$ws.Cells.Item(127, 96).Value = 56
$ws.Cells.Item(127, 97).Value = 16.9
$ws.Cells.Item(127, 98).Value = 25.9
$ws.Cells.Item(127, 99).Value = 71.7
$ws.Cells.Item(127, 100).Value = 3.8
$ws.Cells.Item(127, 101).Value = 7.6
$ws.Cells.Item(127, 102).Value = 0.913
$ws.Cells.Item(127, 103).Value = 0.562
$ws.Cells.Item(127, 104).Value = 17.9
$ws.Cells.Item(127, 105).Value = 24.6
$ws.Cells.Item(127, 106).Value = 2.3
$ws.Cells.Item(127, 107).Value = 0.957
$ws.Cells.Item(127, 108).Value = 0.924
$ws.Cells.Item(127, 109).Value = 0.8
$ws.Cells.Item(127, 110).Value = 0.2
$ws.Cells.Item(127, 111).Value = "'FALSE"
